# Applies the scheduled-runner profit recalculation update to all 8 sheets.
# Values are taken verbatim from the authoritative change set derived from the
# canonical-OOXML diff; columns H-N hold currentAveragePrice / Leve price / profit
# metrics that were recomputed upstream (no formulas live in these cells).

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 356.33334
$ws.Range("I9").Value = 374.625
$ws.Range("J9").Value = 210
$ws.Range("K9").Value = 374.625
$ws.Range("L9").Value = 210
$ws.Range("M9").Value = -205.625
$ws.Range("N9").Value = -548
$ws.Range("H15").Value = 1348.4255
$ws.Range("I15").Value = 1348.4255
$ws.Range("K15").Value = 4045.2765
$ws.Range("M15").Value = -3876.2765
$ws.Range("H40").Value = 1583.909
$ws.Range("I40").Value = 1479.5
$ws.Range("K40").Value = 1479.5
$ws.Range("M40").Value = -1304.5
$ws.Range("H44").Value = 38000
$ws.Range("J44").Value = 38000
$ws.Range("L44").Value = 38000
$ws.Range("N44").Value = -38924
$ws.Range("H62").Value = 2062.5
$ws.Range("I62").Value = 2062.5
$ws.Range("K62").Value = 2062.5
$ws.Range("M62").Value = -1438.5
$ws.Range("H65").Value = 2062.5
$ws.Range("I65").Value = 2062.5
$ws.Range("K65").Value = 10312.5
$ws.Range("M65").Value = -7192.5
$ws.Range("H69").Value = 15020
$ws.Range("I69").Value = 13360
$ws.Range("K69").Value = 40080
$ws.Range("M69").Value = -39206
$ws.Range("H70").Value = 1490.4286
$ws.Range("J70").Value = 1480
$ws.Range("L70").Value = 4440
$ws.Range("N70").Value = -4980
$ws.Range("H72").Value = 15020
$ws.Range("I72").Value = 13360
$ws.Range("K72").Value = 120240
$ws.Range("M72").Value = -115872
$ws.Range("H73").Value = 1490.4286
$ws.Range("J73").Value = 1480
$ws.Range("L73").Value = 4440
$ws.Range("N73").Value = -6312
$ws.Range("H88").Value = 2838.318
$ws.Range("I88").Value = 4406.25
$ws.Range("J88").Value = 1942.3572
$ws.Range("K88").Value = 4406.25
$ws.Range("L88").Value = 1942.3572
$ws.Range("M88").Value = -4000.25
$ws.Range("N88").Value = -2754.3572
$ws.Range("H91").Value = 2838.318
$ws.Range("I91").Value = 4406.25
$ws.Range("J91").Value = 1942.3572
$ws.Range("K91").Value = 4406.25
$ws.Range("L91").Value = 1942.3572
$ws.Range("M91").Value = -3002.25
$ws.Range("N91").Value = -4750.3572
$ws.Range("H112").Value = 3231.1667
$ws.Range("I112").Value = 2397
$ws.Range("J112").Value = 4399
$ws.Range("K112").Value = 7191
$ws.Range("L112").Value = 13197
$ws.Range("M112").Value = -6083
$ws.Range("N112").Value = -15413
$ws.Range("H135").Value = 3739.5625
$ws.Range("I135").Value = 3630.9285
$ws.Range("K135").Value = 32678.3565
$ws.Range("M135").Value = -30143.3565
$ws.Range("H137").Value = 18522254
$ws.Range("I137").Value = 26317908
$ws.Range("J137").Value = 7574.25
$ws.Range("K137").Value = 78953724
$ws.Range("L137").Value = 22722.75
$ws.Range("M137").Value = -78951174
$ws.Range("N137").Value = -27822.75
$ws.Range("H138").Value = 2165.5
$ws.Range("I138").Value = 1122.375
$ws.Range("K138").Value = 3367.125
$ws.Range("M138").Value = 1772.875
$ws.Range("H141").Value = 1778.4
$ws.Range("I141").Value = 1590.8148
$ws.Range("K141").Value = 4772.4444
$ws.Range("M141").Value = 407.5555999999997

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 868728
$ws.Range("I32").Value = 1068013.4
$ws.Range("K32").Value = 1068013.4
$ws.Range("M32").Value = -1067726.4
$ws.Range("H45").Value = 5252
$ws.Range("I45").Value = 4837.3335
$ws.Range("J45").Value = 5666.6665
$ws.Range("K45").Value = 4837.3335
$ws.Range("L45").Value = 5666.6665
$ws.Range("M45").Value = -4460.3335
$ws.Range("N45").Value = -6420.6665

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1565.4865
$ws.Range("I86").Value = 1341.5
$ws.Range("J86").Value = 2525.4285
$ws.Range("K86").Value = 1341.5
$ws.Range("L86").Value = 2525.4285
$ws.Range("M86").Value = -218.5
$ws.Range("N86").Value = -4771.4285
$ws.Range("H89").Value = 1565.4865
$ws.Range("I89").Value = 1341.5
$ws.Range("J89").Value = 2525.4285
$ws.Range("K89").Value = 6707.5
$ws.Range("L89").Value = 12627.1425
$ws.Range("M89").Value = -1091.5
$ws.Range("N89").Value = -23859.1425

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("N21").ClearContents()
$ws.Range("H22").Value = 1289.3846
$ws.Range("I22").Value = 156
$ws.Range("J22").Value = 2611.6667
$ws.Range("K22").Value = 156
$ws.Range("L22").Value = 2611.6667
$ws.Range("M22").Value = 194
$ws.Range("N22").Value = -3311.6667
$ws.Range("H31").Value = 1237477.1
$ws.Range("I31").Value = 1464757.5
$ws.Range("J31").Value = 3669.5715
$ws.Range("K31").Value = 1464757.5
$ws.Range("L31").Value = 3669.5715
$ws.Range("M31").Value = -1464462.5
$ws.Range("N31").Value = -4259.5715
$ws.Range("H34").Value = 1237477.1
$ws.Range("I34").Value = 1464757.5
$ws.Range("J34").Value = 3669.5715
$ws.Range("K34").Value = 1464757.5
$ws.Range("L34").Value = 3669.5715
$ws.Range("M34").Value = -1464555.5
$ws.Range("N34").Value = -4073.5715
$ws.Range("H39").Value = 2033.3334
$ws.Range("J39").Value = 2575
$ws.Range("L39").Value = 2575
$ws.Range("N39").Value = -3357
$ws.Range("H49").Value = 2033.3334
$ws.Range("J49").Value = 2575
$ws.Range("L49").Value = 2575
$ws.Range("N49").Value = -2939
$ws.Range("H132").Value = 4473.143
$ws.Range("I132").Value = 3763.2727
$ws.Range("J132").Value = 7076
$ws.Range("K132").Value = 11289.8181
$ws.Range("L132").Value = 21228
$ws.Range("M132").Value = -8759.8181
$ws.Range("N132").Value = -26288
$ws.Range("H134").Value = 2589.1746
$ws.Range("I134").Value = 2555.5264
$ws.Range("J134").Value = 2640.32
$ws.Range("K134").Value = 7666.5792
$ws.Range("L134").Value = 7920.960000000001
$ws.Range("M134").Value = -5131.5792
$ws.Range("N134").Value = -12990.96

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2269627.8
$ws.Range("I5").Value = 2551966.8
$ws.Range("K5").Value = 7655900.399999999
$ws.Range("M5").Value = -7655788.399999999
$ws.Range("H23").Value = 401.92
$ws.Range("I23").Value = 255.09091
$ws.Range("J23").Value = 517.2857
$ws.Range("K23").Value = 765.27273
$ws.Range("L23").Value = 1551.8571
$ws.Range("M23").Value = -530.27273
$ws.Range("N23").Value = -2021.8571
$ws.Range("H34").Value = 3471.111
$ws.Range("J34").Value = 4833.3335
$ws.Range("L34").Value = 14500.0005
$ws.Range("N34").Value = -14668.0005
$ws.Range("H38").Value = 73.13333
$ws.Range("I38").Value = 162.25
$ws.Range("K38").Value = 486.75
$ws.Range("M38").Value = -139.75
$ws.Range("H55").Value = 105003270
$ws.Range("J55").Value = 6275
$ws.Range("L55").Value = 18825
$ws.Range("N55").Value = -19179
$ws.Range("H92").Value = 822.4138
$ws.Range("I92").Value = 523.8095
$ws.Range("J92").Value = 1606.25
$ws.Range("K92").Value = 1571.4285
$ws.Range("L92").Value = 4818.75
$ws.Range("M92").Value = -323.4285
$ws.Range("N92").Value = -7314.75
$ws.Range("H109").Value = 2609.889
$ws.Range("I109").Value = 1498
$ws.Range("J109").Value = 3999.75
$ws.Range("K109").Value = 4494
$ws.Range("L109").Value = 11999.25
$ws.Range("M109").Value = -3454
$ws.Range("N109").Value = -14079.25
$ws.Range("H114").Value = 4065.4
$ws.Range("J114").Value = 4306
$ws.Range("L114").Value = 12918
$ws.Range("N114").Value = -19426
$ws.Range("H135").Value = 2269627.8
$ws.Range("I135").Value = 2551966.8
$ws.Range("K135").Value = 22967701.2
$ws.Range("M135").Value = -22965166.2
$ws.Range("H137").Value = 3184.5454
$ws.Range("I137").Value = 2055
$ws.Range("K137").Value = 6165
$ws.Range("M137").Value = -1065

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H126").Value = 5056.8887
$ws.Range("I126").Value = 5003.5835
$ws.Range("K126").Value = 15010.7505
$ws.Range("M126").Value = -12540.7505
$ws.Range("H135").Value = 78700.53999999999
$ws.Range("J135").Value = 78700.53999999999
$ws.Range("L135").Value = 78700.53999999999
$ws.Range("N135").Value = -88840.53999999999

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2262.16
$ws.Range("I40").Value = 2191
$ws.Range("J40").Value = 2388.6667
$ws.Range("K40").Value = 2191
$ws.Range("L40").Value = 2388.6667
$ws.Range("M40").Value = -2055
$ws.Range("N40").Value = -2660.6667
$ws.Range("H122").Value = 5801
$ws.Range("I122").Value = 4666.6665
$ws.Range("K122").Value = 13999.9995
$ws.Range("M122").Value = -11549.9995

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 24428.2
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
